$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
